$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.161.27"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "1.835.67"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'242.05"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").Value = "'0.6620"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.07414"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("D9").Value = "'0.2941"
$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").Value = "'23.02"
$ws.Range("E10").Value = "  +1.62%  "

$ws.Range("D11").Value = "'0.07751"
$ws.Range("E11").Value = "  +1.15%  "

$ws.Range("D12").Value = "1.840.31"
$ws.Range("E12").Value = "  +4.87%  "

$ws.Range("D13").Value = "'5.003"
$ws.Range("E13").Value = "  -0.07%  "

$ws.Range("D14").Value = "'0.6704"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").Value = "'83.08"
$ws.Range("E15").Value = "  -3.26%  "

$ws.Range("D16").Value = "'6.122"
$ws.Range("E16").Value = "  -0.70%  "

$ws.Range("D17").Value = "'0.000008612"
$ws.Range("E17").Value = "  +5.07%  "

$ws.Range("D18").Value = "29.150.14"
$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("D19").Value = "2.083.12"
$ws.Range("E19").Value = "  +0.25%  "

$ws.Range("D20").Value = "'227.19"
$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("D21").Value = "'12.50"
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("D23").Value = "'7.133"
$ws.Range("E23").Value = "  -1.73%  "

$ws.Range("D24").Value = "'0.9998"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").Value = "'158.96"
$ws.Range("E25").Value = "  -0.97%  "

$ws.Range("D26").Value = "'0.1410"
$ws.Range("E26").Value = "  +0.46%  "

$ws.Range("D27").Value = "'8.605"
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("D28").Value = "'18.02"
$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").Value = "'1.511"

$ws.Range("D30").Value = "'4.125"
$ws.Range("E30").Value = "  -1.88%  "

$ws.Range("D31").Value = "'4.057"
$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("D32").Value = "'1.191"
$ws.Range("E32").Value = "  +1.08%  "

$ws.Range("D33").Value = "'0.05295"
$ws.Range("E33").Value = "  -0.85%  "

$ws.Range("D34").Value = "'1.873"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("D35").Value = "'0.7396"
$ws.Range("E35").Value = "  -0.73%  "

$ws.Range("D36").Value = "'1.148"
$ws.Range("E36").Value = "  +1.72%  "

$ws.Range("D37").Value = "'2.649"
$ws.Range("E37").Value = "  -1.25%  "

$ws.Range("D38").Value = "1.303.15"
$ws.Range("E38").Value = "  -1.00%  "

$ws.Range("D39").Value = "'0.01792"
$ws.Range("E39").Value = "  -0.42%  "

$ws.Range("D40").Value = "'2.739"
$ws.Range("E40").Value = "  +1.05%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.314"
$ws.Range("E41").Value = "  +5.90%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.9174"
$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'0.9995"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("B44").Value = "XinFinNetwork"
$ws.Range("C44").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D44").Value = "'0.08259"
$ws.Range("E44").Value = "  +10.20%  "

$ws.Range("D45").Value = "'102.61"
$ws.Range("E45").Value = "  -0.63%  "

$ws.Range("D46").Value = "1.983.95"
$ws.Range("E46").Value = "  +2.64%  "

$ws.Range("D47").Value = "'0.5136"
$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("D48").Value = "'64.16"
$ws.Range("E48").Value = "  +0.75%  "

$ws.Range("D49").Value = "'1.745"
$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000117"
$ws.Range("E50").Value = "  -4.75%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05842"
$ws.Range("E51").Value = "  -1.31%  "

